$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# These cells already hold numeric-looking values ("11.1", "18.2", ...)
# stored as plain text. A direct $range.Value = "11.13" assignment would
# get auto-coerced to a number by Excel, changing the cell's stored type.
# To keep them as text (matching the source data), enter each value as a
# text-producing formula and then collapse it to a static value via
# copy / paste-special-values, which preserves the Text cell type
# without leaving a residual formula behind.
function Set-TextValue($range, [string]$text) {
    $range.Formula = '="' + $text + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null
}

# Enterprises density (per 1000 people) row
Set-TextValue $ws.Range("B11") "11.13"
Set-TextValue $ws.Range("D11") "11.53"

# Employment (% of total) row
Set-TextValue $ws.Range("C12") "18.23"
Set-TextValue $ws.Range("D12") "77.23"

# Enterprises (% of total) row
Set-TextValue $ws.Range("B14") "96.34"
Set-TextValue $ws.Range("C14") "3.47"
Set-TextValue $ws.Range("D14") "99.81"
